# Apply the edit: insert two new price rows (rows 14 and 15) into the
# "Ciruela" daily price sheet, pushing the existing rows 14-83 down to 16-85.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 14 (Excel's Insert copies the
# formatting - including the date number format on column D - from the row
# above, same as a manual "Insert Copied Cells" / row insert in the UI).
$ws.Rows("14:15").Insert()

# --- New row 14 ---
$ws.Range("A14").Value = 5
$ws.Range("B14").Value = "Macroferia Regional de Talca"
$ws.Range("C14").Value = "Maule"
$ws.Range("D14").Value = 44608
$ws.Range("E14").Value = 7
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100103
$ws.Range("H14").Value = "Frutos de hueso (carozo)"
$ws.Range("I14").Value = 100103002
$ws.Range("J14").Value = "Ciruela"
$ws.Range("K14").Value = "Black Amber"
$ws.Range("L14").Value = "Especial"
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 13000
$ws.Range("O14").Value = 13000
$ws.Range("P14").Value = 13000
$ws.Range("Q14").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R14").Value = "Región de O'Higgins"
$ws.Range("S14").Value = 722
$ws.Range("T14").Value = 18

# --- New row 15 ---
$ws.Range("A15").Value = 5
$ws.Range("B15").Value = "Macroferia Regional de Talca"
$ws.Range("C15").Value = "Maule"
$ws.Range("D15").Value = 44608
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100103
$ws.Range("H15").Value = "Frutos de hueso (carozo)"
$ws.Range("I15").Value = 100103002
$ws.Range("J15").Value = "Ciruela"
$ws.Range("K15").Value = "Black Amber"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 11000
$ws.Range("O15").Value = 11000
$ws.Range("P15").Value = 11000
$ws.Range("Q15").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R15").Value = "Región de O'Higgins"
$ws.Range("S15").Value = 611
$ws.Range("T15").Value = 18
